$wb = $excel.ActiveWorkbook

# The localization status for the two entries moved from "Ready for handoff"
# to "In Translation" - update every sheet that surfaces the Status column.

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"

# Re-generating the report re-sizes the (now shorter) status columns to fit
# the new text.
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5
$ws2.Columns.Item(3).ColumnWidth = 12.5
$ws3.Columns.Item(3).ColumnWidth = 12.5
